$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Crafting quest instructions (row 6, column D): add a reminder paragraph ---
$oldCraft = @'
<p><strong>Crafting cannot be automated. Keep an eye on Server Messages section to see successes, failures and if you have new items to craft.</strong></p><p><strong>Desktop:</strong></p>
'@
$newCraft = @'
<p><strong>Crafting cannot be automated. Keep an eye on Server Messages section to see successes, failures and if you have new items to craft.</strong></p><p><strong><em>Remember to keep crafting - even after this quest. Stick to weapons (any type) - but do experiment with others, like Armour. A later quest will be less painful.</em></strong></p><p><strong>Desktop:</strong></p>
'@
$craftText = $ws.Range("D6").Value2
$ws.Range("D6").Value2 = $craftText.Replace($oldCraft, $newCraft)

# --- 2) Enchanting quest instructions (row 9, column D): add a reminder paragraph ---
$oldEnchant = @'
So lets get a lot of it.</p><p><strong>Desktop/Mobile:</strong></p>
'@
$newEnchant = @'
So lets get a lot of it.</p><p><strong><em>Remember to keep enchanting and disenchanting items. A later quest will be less painful.</em></strong></p><p><strong>Desktop/Mobile:</strong></p>
'@
$enchantText = $ws.Range("D9").Value2
$ws.Range("D9").Value2 = $enchantText.Replace($oldEnchant, $newEnchant)

# --- 3) Class skills quest instructions (row 11, column D): remove the Heretics example paragraphs ---
$oldCS = @'
- Attack type</p><p>\For example, Heretics:</p><p>With a damage spell equipped you have a small chance to cast another spell. Enemies cannot avoid this.</p><p>This means while casting and with at least one damage spell, based on the % of your class bonus you can cast another spell, automatically. Each class has its own special which you can read about in the help docs for your specific class.</p><p>To do this:</p>
'@
$newCS = @'
- Attack type</p><p>To do this:</p>
'@
$csText = $ws.Range("D11").Value2
$ws.Range("D11").Value2 = $csText.Replace($oldCS, $newCS)

# --- 4) Mercenary for hire quest intro text (row 12, column C): fix "poet"/"Port" typos ---
$oldMerc = @'
Where did they go? You ask The poet in a rush, frantic and worried about your parents.<br /> <br /> “Child. Calm your self. They are here, some where. You’ll have to find them eventually.”<br /> <br /> Find them? What has he done with them? You start to get angry but The Port waves his hand and calmness washes over you.
'@
$newMerc = @'
Where did they go? You ask The Poet in a rush, frantic and worried about your parents.<br /> <br /> “Child. Calm your self. They are here, somewhere. You’ll have to find them eventually.”<br /> <br /> Find them? What has he done with them? You start to get angry but The Poet waves his hand and calmness washes over you.
'@
$mercText = $ws.Range("C12").Value2
$ws.Range("C12").Value2 = $mercText.Replace($oldMerc, $newMerc)

# --- 5) Row 10 (id 9): remove I10 value and change AE10 500 -> 350 ---
$ws.Range("I10").ClearContents()
$ws.Range("AE10").Value = 350

# --- 6) Row 11 (id 10): E11 350 -> 250, K11 25 -> 35, AE11 600 -> 450 ---
$ws.Range("E11").Value = 250
$ws.Range("K11").Value = 35
$ws.Range("AE11").Value = 450

# --- 7) Row 12 (id 11): add required_quest_id (T12) referencing new quest name ---
$ws.Range("T12").Value = "The truth is out there"
